$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New commit rows: text in column C (using the same highlighted cell style
# already used by the rows above) and hours in column G.
$ws.Range("C10").Value = "move & look instruction added"
$ws.Range("C10").Style = "20% - Énfasis5"
$ws.Range("G10").Value = 2.3

$ws.Range("C11").Value = "pick instruction added"
$ws.Range("C11").Style = "20% - Énfasis5"
$ws.Range("G11").Value = 2

$ws.Range("C12").Value = "look update"
$ws.Range("C12").Style = "20% - Énfasis5"
$ws.Range("G12").Value = 0.3

$ws.Range("C13").Value = "object added to the game data"
$ws.Range("C13").Style = "20% - Énfasis5"
$ws.Range("G13").Value = 0.4

$ws.Range("C14").Value = "general update (class TYPES) "
$ws.Range("C14").Style = "20% - Énfasis5"
$ws.Range("G14").Value = 0.5

# Update the selected cell to match the saved view state.
$ws.Range("C16").Select()
